# Add a new "Length_fov(mm)" column between the existing "Global_Uncertainty"
# column and the "RealLength(cm)" column, shifting "RealLength(cm)" and
# "Pond_Type" one column to the right (T->U, U->V), and recompute
# "RealLength(cm)" using an updated conversion factor that accounts for the
# diagonal of the carapace bounding box (per commit message "add diagonal to
# carapace").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank column at T (column 20); everything from T onward
# shifts right by one (old T -> U, old U -> V).
$ws.Columns.Item(20).Insert()

# New header for the inserted column.
$ws.Range("T1").Value = "Length_fov(mm)"

# Conversion factors derived from the committed data: the previously
# computed "RealLength(cm)" values (now sitting in column U after the
# insert) are the basis for both the new "Length_fov(mm)" column (T) and
# the recomputed "RealLength(cm)" column (U).
$fovFactor = 0.9761170231700621
$realLengthFactor = 0.9797734627831713

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $oldCell = $ws.Cells.Item($r, 21)   # column U now holds the old RealLength(cm) value
    $oldValue = $oldCell.Value()

    if ($oldValue -eq "") {
        # Row had no measurement (original value was NaN) - leave both
        # Length_fov(mm) and RealLength(cm) blank.
        continue
    }

    $numericValue = [double]$oldValue

    $fovCell = $ws.Cells.Item($r, 20)   # column T - new Length_fov(mm)
    $fovCell.Value = $numericValue * $fovFactor

    $oldCell.Value = $numericValue * $realLengthFactor
}
